$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values (B5:AH5) to 2 decimal places, matching target data
$row5 = @{
    "B5" = 14.96;  "C5" = 10.92;  "D5" = 0.7;    "E5" = 31.69;  "F5" = 25.72;
    "G5" = 11.01;  "H5" = 44.02;  "I5" = 17.21;  "J5" = 7.6;    "K5" = 11.42;
    "L5" = 12.58;  "M5" = 13.84;  "N5" = 3.62;   "O5" = 10.86;  "P5" = 16.11;
    "Q5" = 9.31;   "R5" = 0.42;   "S5" = 0.33;   "T5" = 164.17; "U5" = 31.53;
    "V5" = 10.6;   "W5" = 21.5;   "X5" = 11.37;  "Y5" = 1.48;   "Z5" = 21.27;
    "AA5" = 9.19;  "AB5" = 7.79;  "AC5" = 9.58;  "AD5" = 13.34; "AE5" = 0.12;
    "AF5" = 39.46; "AG5" = 5.78;  "AH5" = 12.84
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Delete row 6 entirely
$ws.Rows.Item(6).Delete()
